# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data between row 211 and row 212 ---
# (id/A, Div/C, Div Original Name/D and Date/E stay attached to their own
# row; everything else - the match id in B and all the odds data in F:AC -
# belongs to the other match and needs to move to the other row.)

# Save row 212's current values (this is what row 211 needs to become)
$b212 = $ws.Range("B212").Value()
$f212 = $ws.Range("F212").Value()
$g212 = $ws.Range("G212").Value()
$h212 = $ws.Range("H212").Value()
$i212 = $ws.Range("I212").Value()
$j212 = $ws.Range("J212").Value()
$k212 = $ws.Range("K212").Value()
$l212 = $ws.Range("L212").Value()
$m212 = $ws.Range("M212").Value()
$n212 = $ws.Range("N212").Value()
$o212 = $ws.Range("O212").Value()
$p212 = $ws.Range("P212").Value()
$q212 = $ws.Range("Q212").Value()
$r212 = $ws.Range("R212").Value()
$s212 = $ws.Range("S212").Value()
$t212 = $ws.Range("T212").Value()
$u212 = $ws.Range("U212").Value()
$v212 = $ws.Range("V212").Value()
$w212 = $ws.Range("W212").Value()
$x212 = $ws.Range("X212").Value()
$y212 = $ws.Range("Y212").Value()
$z212 = $ws.Range("Z212").Value()
$aa212 = $ws.Range("AA212").Value()
$ab212 = $ws.Range("AB212").Value()
$ac212 = $ws.Range("AC212").Value()

# Save row 211's current values (this is what row 212 needs to become)
$b211 = $ws.Range("B211").Value()
$f211 = $ws.Range("F211").Value()
$g211 = $ws.Range("G211").Value()
$h211 = $ws.Range("H211").Value()
$i211 = $ws.Range("I211").Value()
$j211 = $ws.Range("J211").Value()
$k211 = $ws.Range("K211").Value()
$l211 = $ws.Range("L211").Value()
$m211 = $ws.Range("M211").Value()
$n211 = $ws.Range("N211").Value()
$o211 = $ws.Range("O211").Value()
$p211 = $ws.Range("P211").Value()
$q211 = $ws.Range("Q211").Value()
$r211 = $ws.Range("R211").Value()
$s211 = $ws.Range("S211").Value()
$t211 = $ws.Range("T211").Value()
$u211 = $ws.Range("U211").Value()
$v211 = $ws.Range("V211").Value()
$w211 = $ws.Range("W211").Value()
$x211 = $ws.Range("X211").Value()
$y211 = $ws.Range("Y211").Value()
$z211 = $ws.Range("Z211").Value()
$aa211 = $ws.Range("AA211").Value()
$ab211 = $ws.Range("AB211").Value()
$ac211 = $ws.Range("AC211").Value()

# Write row 212's former values into row 211
$ws.Range("B211").Value = $b212
$ws.Range("F211").Value = $f212
$ws.Range("G211").Value = $g212
$ws.Range("H211").Value = $h212
$ws.Range("I211").Value = $i212
$ws.Range("J211").Value = $j212
$ws.Range("K211").Value = $k212
$ws.Range("L211").Value = $l212
$ws.Range("M211").Value = $m212
$ws.Range("N211").Value = $n212
$ws.Range("O211").Value = $o212
$ws.Range("P211").Value = $p212
$ws.Range("Q211").Value = $q212
$ws.Range("R211").Value = $r212
$ws.Range("S211").Value = $s212
$ws.Range("T211").Value = $t212
$ws.Range("U211").Value = $u212
$ws.Range("V211").Value = $v212
$ws.Range("W211").Value = $w212
$ws.Range("X211").Value = $x212
$ws.Range("Y211").Value = $y212
$ws.Range("Z211").Value = $z212
$ws.Range("AA211").Value = $aa212
$ws.Range("AB211").Value = $ab212
$ws.Range("AC211").Value = $ac212

# Write row 211's former values into row 212
$ws.Range("B212").Value = $b211
$ws.Range("F212").Value = $f211
$ws.Range("G212").Value = $g211
$ws.Range("H212").Value = $h211
$ws.Range("I212").Value = $i211
$ws.Range("J212").Value = $j211
$ws.Range("K212").Value = $k211
$ws.Range("L212").Value = $l211
$ws.Range("M212").Value = $m211
$ws.Range("N212").Value = $n211
$ws.Range("O212").Value = $o211
$ws.Range("P212").Value = $p211
$ws.Range("Q212").Value = $q211
$ws.Range("R212").Value = $r211
$ws.Range("S212").Value = $s211
$ws.Range("T212").Value = $t211
$ws.Range("U212").Value = $u211
$ws.Range("V212").Value = $v211
$ws.Range("W212").Value = $w211
$ws.Range("X212").Value = $x211
$ws.Range("Y212").Value = $y211
$ws.Range("Z212").Value = $z211
$ws.Range("AA212").Value = $aa211
$ws.Range("AB212").Value = $ab211
$ws.Range("AC212").Value = $ac211

# --- Append a new match as row 218 ---
# Copy row 217's cell formats down into row 218 first, so the new row
# picks up the same styling (bold/bordered id cell, date number format)
# used throughout the table.
$ws.Range("A217:AC217").Copy() | Out-Null
$ws.Range("A218:AC218").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A218").Value = 216
$ws.Range("B218").Value = 8021846
$ws.Range("C218").Value = "Bolivia Primera División"
$ws.Range("D218").Value = "Bolivia Apertura"
$ws.Range("E218").Value = 45393.83333333334
$ws.Range("F218").Value = "Real Santa Cruz"
$ws.Range("G218").Value = "The Strongest"
$ws.Range("K218").Value = 6
$ws.Range("L218").Value = 5
$ws.Range("M218").Value = 1.333
$ws.Range("N218").Value = 6.5
$ws.Range("O218").Value = 5
$ws.Range("P218").Value = 1.444
$ws.Range("Q218").Value = 1
$ws.Range("R218").Value = 2.05
$ws.Range("S218").Value = 1.75
$ws.Range("T218").Value = 3
$ws.Range("U218").Value = 2.025
$ws.Range("V218").Value = 1.775
$ws.Range("W218").Value = 0
$ws.Range("X218").Value = 0
$ws.Range("Y218").Value = 0
$ws.Range("Z218").Value = 0
$ws.Range("AA218").Value = 0

Write-Output "Edit complete"
